$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.403.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.907.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.904.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.442"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.121"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.388.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.908.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.436.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "413.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0960"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.922"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0673"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "371.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0341"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.680.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.234"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  -0.79%  "
